# Update a handful of cell values in the Skills Audit worksheet to fix
# wording / typos, as described by the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D13").Value = "Course from Leeds Beckett University website and youtube videos "
$ws.Range("E13").Value = "Get more knowledge to create eerd, erd and uml diagrams."
$ws.Range("E14").Value = "Obtain more knowledge to make documentation better."
$ws.Range("E15").Value = "Gain more practice and ideas to make a better logo."
$ws.Range("E17").Value = "Obtain proper information to make the website reponsive."
